# Commit "#5: fund, bonds, otherbonds, antique done"
#
# The "基金受益憑證" (fund) sheet is dropped. The data that used to live
# half in that placeholder sheet and half (mis-shaped, without a
# "quantity"/"otherbonds" column) in "具有相當價值之財產" (other valuable
# property) is consolidated into a single, properly-shaped table that now
# lives in "具有相當價值之財產": name / quantity / owner / total /
# property_category / category / date / legislator_name / legislator_id /
# source_file / index.
#
# All of the sheets after it (保險/債務/事業投資) are untouched content-wise;
# Excel simply slides them up one slot once "基金受益憑證" disappears.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a text value into $range without Excel auto-coercing a
# "yyyy-mm-dd"-shaped string into a date serial. We stage the literal
# value in a scratch cell that has been force-formatted as Text, copy
# it, and paste-special (values only) onto the destination -- the
# destination keeps its own (default/General) style while inheriting
# the already-stringified value.
# ---------------------------------------------------------------------
function Set-TextValue($ws, $range, [string]$text) {
    $scratch = $ws.Range("ZZ500")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# 1. Drop the "基金受益憑證" sheet entirely.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("基金受益憑證").Delete()

# ---------------------------------------------------------------------
# 2. Rebuild "具有相當價值之財產" with the new combined dataset.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("具有相當價值之財產")
$ws.Cells.Clear()

# Header row
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "quantity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "property_category"
$ws.Range("G1").Value = "category"
Set-TextValue $ws ($ws.Range("H1")) "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"
$ws.Range("K1").Value = "source_file"
$ws.Range("L1").Value = "index"

# Row 2: 油畫 (oil painting)
$ws.Range("A2").Value = 69
$ws.Range("B2").Value = "油畫"
$ws.Range("C2").Value = "2件"
$ws.Range("D2").Value = "姚文智"
$ws.Range("E2").Value = 400000
$ws.Range("F2").Value = "otherbonds"
$ws.Range("G2").Value = "normal"
Set-TextValue $ws ($ws.Range("H2")) "2012-04-30"
$ws.Range("I2").Value = "姚文智"
$ws.Range("J2").Value = 1745
Set-TextValue $ws ($ws.Range("K2")) "tmp60da1"
$ws.Range("L2").Value = 69

# Row 3: 琉璃 (glazed glass)
$ws.Range("A3").Value = 70
$ws.Range("B3").Value = "琉璃"
$ws.Range("C3").Value = "2件"
$ws.Range("D3").Value = "姚文智"
$ws.Range("E3").Value = 250000
$ws.Range("F3").Value = "otherbonds"
$ws.Range("G3").Value = "normal"
Set-TextValue $ws ($ws.Range("H3")) "2012-04-30"
$ws.Range("I3").Value = "姚文智"
$ws.Range("J3").Value = 1745
Set-TextValue $ws ($ws.Range("K3")) "tmp60da1"
$ws.Range("L3").Value = 70

# Row 4: 珠寶 (jewelry)
$ws.Range("A4").Value = 71
$ws.Range("B4").Value = "珠寶"
$ws.Range("C4").Value = "3件"
$ws.Range("D4").Value = "潘瓊琪"
$ws.Range("E4").Value = 500000
$ws.Range("F4").Value = "otherbonds"
$ws.Range("G4").Value = "normal"
Set-TextValue $ws ($ws.Range("H4")) "2012-04-30"
$ws.Range("I4").Value = "姚文智"
$ws.Range("J4").Value = 1745
Set-TextValue $ws ($ws.Range("K4")) "tmp60da1"
$ws.Range("L4").Value = 71

# ---------------------------------------------------------------------
# 3. Match the workbook-wide formatting convention: header row (row 1)
#    and the whole index column (A) are bold + thin-bordered + centered
#    / top-aligned; every other data cell keeps the plain default style.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("B1:L1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$indexRange = $ws.Range("A2:A4")
$indexRange.Font.Bold = $true
$indexRange.Borders.LineStyle = 1
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
